$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2546499
$ws.Range("J17").Value = 2546499
$ws.Range("L17").Value = 7639497
$ws.Range("N17").Value = -7639833

$ws.Range("H98").Value = 4267.5835
$ws.Range("I98").Value = 3121.1
$ws.Range("J98").Value = 10000
$ws.Range("K98").Value = 3121.1
$ws.Range("L98").Value = 10000
$ws.Range("M98").Value = -1623.1
$ws.Range("N98").Value = -12996

$ws.Range("H112").Value = 4879.75
$ws.Range("I112").Value = 745
$ws.Range("J112").Value = 5032.8887
$ws.Range("K112").Value = 2235
$ws.Range("L112").Value = 15098.6661
$ws.Range("M112").Value = -1127
$ws.Range("N112").Value = -17314.6661

$ws.Range("H122").Value = 4267.5835
$ws.Range("I122").Value = 3121.1
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 9363.299999999999
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -6913.299999999999
$ws.Range("N122").Value = -34900

$ws.Range("H137").Value = 1503.2727
$ws.Range("I137").Value = 1230.6666
$ws.Range("J137").Value = 2730
$ws.Range("K137").Value = 3691.9998
$ws.Range("L137").Value = 8190
$ws.Range("M137").Value = -1141.9998
$ws.Range("N137").Value = -13290

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3352.2
$ws.Range("I2").Value = 3190.125
$ws.Range("K2").Value = 3190.125
$ws.Range("M2").Value = -3077.125

$ws.Range("H5").Value = 98.75
$ws.Range("I5").Value = 98.333336
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 98.333336
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 13.666664
$ws.Range("N5").Value = -324

$ws.Range("H45").Value = 5003.077
$ws.Range("I45").Value = 6935.8823
$ws.Range("K45").Value = 6935.8823
$ws.Range("M45").Value = -6558.8823

$ws.Range("H61").Value = 6826.7144
$ws.Range("I61").Value = 7340.0527
$ws.Range("K61").Value = 7340.0527
$ws.Range("M61").Value = -7128.0527

$ws.Range("H74").Value = 1792
$ws.Range("I74").Value = 1705.098
$ws.Range("K74").Value = 1705.098
$ws.Range("M74").Value = -831.098

$ws.Range("H77").Value = 1792
$ws.Range("I77").Value = 1705.098
$ws.Range("K77").Value = 8525.49
$ws.Range("M77").Value = -4157.49

$ws.Range("H97").Value = 1849
$ws.Range("I97").Value = 2282.2856
$ws.Range("J97").Value = 1512
$ws.Range("K97").Value = 2282.2856
$ws.Range("L97").Value = 1512
$ws.Range("M97").Value = -1786.2856
$ws.Range("N97").Value = -2504

$ws.Range("H116").Value = 3352.2
$ws.Range("I116").Value = 3190.125
$ws.Range("K116").Value = 3190.125
$ws.Range("M116").Value = -896.125

$ws.Range("H122").Value = 1285248.9
$ws.Range("I122").Value = 1604871
$ws.Range("K122").Value = 4814613
$ws.Range("M122").Value = -4812163

$ws.Range("H132").Value = 3443.425
$ws.Range("I132").Value = 1748.6
$ws.Range("J132").Value = 5138.25
$ws.Range("K132").Value = 5245.799999999999
$ws.Range("L132").Value = 15414.75
$ws.Range("M132").Value = -2715.799999999999
$ws.Range("N132").Value = -20474.75

$ws.Range("H136").Value = 6826.7144
$ws.Range("I136").Value = 7340.0527
$ws.Range("K136").Value = 22020.1581
$ws.Range("M136").Value = -19470.1581

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3352.2
$ws.Range("I3").Value = 3190.125
$ws.Range("K3").Value = 3190.125
$ws.Range("M3").Value = -3076.125

$ws.Range("H4").Value = 98.75
$ws.Range("I4").Value = 98.333336
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 98.333336
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = 16.666664
$ws.Range("N4").Value = -330

$ws.Range("H19").Value = 12000
$ws.Range("J19").Value = 12000
$ws.Range("L19").Value = 12000
$ws.Range("N19").Value = -12346

$ws.Range("H99").Value = 66667860
$ws.Range("I99").Value = 90910216
$ws.Range("J99").Value = 1377.75
$ws.Range("K99").Value = 90910216
$ws.Range("L99").Value = 1377.75
$ws.Range("M99").Value = -90908718
$ws.Range("N99").Value = -4373.75

$ws.Range("H105").Value = 16723.857
$ws.Range("I105").Value = 21609
$ws.Range("J105").Value = 4511
$ws.Range("K105").Value = 21609
$ws.Range("L105").Value = 4511
$ws.Range("M105").Value = -19862
$ws.Range("N105").Value = -8005

$ws.Range("H132").Value = 42186.668
$ws.Range("J132").Value = 42186.668
$ws.Range("L132").Value = 42186.668
$ws.Range("N132").Value = -52306.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5918381.5
$ws.Range("I16").Value = 9616382
$ws.Range("K16").Value = 9616382
$ws.Range("M16").Value = -9616095

$ws.Range("H31").Value = 6569.3706
$ws.Range("I31").Value = 1653.0555
$ws.Range("J31").Value = 16402
$ws.Range("K31").Value = 1653.0555
$ws.Range("L31").Value = 16402
$ws.Range("M31").Value = -1358.0555
$ws.Range("N31").Value = -16992

$ws.Range("H34").Value = 6569.3706
$ws.Range("I34").Value = 1653.0555
$ws.Range("J34").Value = 16402
$ws.Range("K34").Value = 1653.0555
$ws.Range("L34").Value = 16402
$ws.Range("M34").Value = -1451.0555
$ws.Range("N34").Value = -16806

$ws.Range("H58").Value = 1525.119
$ws.Range("I58").Value = 757.6957
$ws.Range("K58").Value = 757.6957
$ws.Range("M58").Value = -554.6957

$ws.Range("H74").Value = 20774.625
$ws.Range("J74").Value = 20774.625
$ws.Range("L74").Value = 20774.625
$ws.Range("N74").Value = -22522.625

$ws.Range("H77").Value = 20774.625
$ws.Range("J77").Value = 20774.625
$ws.Range("L77").Value = 62323.875
$ws.Range("N77").Value = -71059.875

$ws.Range("H113").Value = 5918381.5
$ws.Range("I113").Value = 9616382
$ws.Range("K113").Value = 9616382
$ws.Range("M113").Value = -9614212

$ws.Range("H122").Value = 1215.9565
$ws.Range("I122").Value = 1059.4445
$ws.Range("K122").Value = 3178.3335
$ws.Range("M122").Value = -728.3335000000002

$ws.Range("H132").Value = 2654.25
$ws.Range("I132").Value = 2219.1333
$ws.Range("K132").Value = 6657.3999
$ws.Range("M132").Value = -4127.3999

$ws.Range("H134").Value = 4303.5
$ws.Range("I134").Value = 5930.4546
$ws.Range("K134").Value = 17791.3638
$ws.Range("M134").Value = -15256.3638

$ws.Range("H136").Value = 1525.119
$ws.Range("I136").Value = 757.6957
$ws.Range("K136").Value = 2273.0871
$ws.Range("M136").Value = 276.9129000000003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4975
$ws.Range("J80").Value = 4975
$ws.Range("L80").Value = 14925
$ws.Range("N80").Value = -16797

$ws.Range("H83").Value = 4975
$ws.Range("J83").Value = 4975
$ws.Range("L83").Value = 44775
$ws.Range("N83").Value = -54135

$ws.Range("H133").Value = 53455.953
$ws.Range("I133").Value = 204406.2
$ws.Range("J133").Value = 9058.823
$ws.Range("K133").Value = 613218.6000000001
$ws.Range("L133").Value = 27176.469
$ws.Range("M133").Value = -608158.6000000001
$ws.Range("N133").Value = -37296.469

$ws.Range("H139").Value = 4470.6523
$ws.Range("I139").Value = 8170
$ws.Range("J139").Value = 2852.1875
$ws.Range("K139").Value = 24510
$ws.Range("L139").Value = 8556.5625
$ws.Range("M139").Value = -19370
$ws.Range("N139").Value = -18836.5625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5080.316
$ws.Range("I132").Value = 9383.333000000001
$ws.Range("J132").Value = 3094.3076
$ws.Range("K132").Value = 28149.999
$ws.Range("L132").Value = 9282.9228
$ws.Range("M132").Value = -25619.999
$ws.Range("N132").Value = -14342.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 166670690
$ws.Range("I40").Value = 200003870
$ws.Range("J40").Value = 4700
$ws.Range("K40").Value = 200003870
$ws.Range("L40").Value = 4700
$ws.Range("M40").Value = -200003734
$ws.Range("N40").Value = -4972

$ws.Range("H136").Value = 4371.3335
$ws.Range("I136").Value = 4818.7188
$ws.Range("J136").Value = 2939.7
$ws.Range("K136").Value = 14456.1564
$ws.Range("L136").Value = 8819.099999999999
$ws.Range("M136").Value = -11906.1564
$ws.Range("N136").Value = -13919.1
